$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D to make room for "Status"; shifts old D..H to E..I
$ws.Columns.Item(4).Insert()

# Remove the row for "Adani Enterprises Limited Rights" (ISIN INE423A20016), which was fully exited
# and no longer appears in the updated dataset; remaining rows shift up by one
$ws.Rows.Item(24).Delete()

# Set header for new Status column
$ws.Range("D1").Value = "Status"

# Update Oct_2025 header (replaces old Nov_2025 in shifted column G)
$ws.Range("G1").Value = "Oct_2025"

# Populate Status + refreshed Oct_2025 / MoM / QoQ figures for each holding row
$ws.Range("D2").Value = "Reducing"
$ws.Range("G2").Value = 6.396431
$ws.Range("H2").Value = -0.04114400000000096
$ws.Range("I2").Value = 3.142735
$ws.Range("D3").Value = "Adding Consistently"
$ws.Range("G3").Value = 9.079192000000001
$ws.Range("H3").Value = 0.4101739999999996
$ws.Range("I3").Value = 0.2216189999999987
$ws.Range("D4").Value = "Adding Consistently"
$ws.Range("G4").Value = 7.112125
$ws.Range("H4").Value = 0.3289010000000001
$ws.Range("I4").Value = 0.6005390000000004
$ws.Range("D5").Value = "Reducing"
$ws.Range("G5").Value = 3.600988
$ws.Range("H5").Value = -0.7395629999999995
$ws.Range("I5").Value = 3.945307
$ws.Range("D6").Value = "Reducing Consistently"
$ws.Range("G6").Value = 8.035861000000001
$ws.Range("H6").Value = -0.5071089999999998
$ws.Range("I6").Value = -0.6257400000000004
$ws.Range("D7").Value = "Reducing"
$ws.Range("G7").Value = 2.972836
$ws.Range("H7").Value = -0.1574119999999999
$ws.Range("I7").Value = 4.190707
$ws.Range("D8").Value = "Adding"
$ws.Range("G8").Value = 6.076973
$ws.Range("H8").Value = 0.1522009999999998
$ws.Range("I8").Value = -0.2832509999999999
$ws.Range("D9").Value = "Adding Consistently"
$ws.Range("G9").Value = 3.882977
$ws.Range("H9").Value = 0.218372
$ws.Range("I9").Value = 0.09590000000000032
$ws.Range("D10").Value = "Adding Consistently"
$ws.Range("G10").Value = 3.272063
$ws.Range("H10").Value = 0.3739309999999998
$ws.Range("I10").Value = 0.5849889999999998
$ws.Range("D11").Value = "Adding"
$ws.Range("G11").Value = 9.659266000000001
$ws.Range("H11").Value = 0.1286550000000002
$ws.Range("I11").Value = -6.775855
$ws.Range("D12").Value = "Reducing Consistently"
$ws.Range("G12").Value = 3.15125
$ws.Range("H12").Value = -0.6120839999999999
$ws.Range("I12").Value = -1.025966
$ws.Range("D13").Value = "Adding Consistently"
$ws.Range("G13").Value = 1.860465
$ws.Range("H13").Value = 0.1122110000000001
$ws.Range("I13").Value = 0.2052830000000001
$ws.Range("D14").Value = "Adding Consistently"
$ws.Range("G14").Value = 1.866245
$ws.Range("H14").Value = 0.0520050000000003
$ws.Range("I14").Value = 0.1859880000000003
$ws.Range("D15").Value = "Reducing Consistently"
$ws.Range("G15").Value = 2.092355
$ws.Range("H15").Value = -0.01621799999999984
$ws.Range("I15").Value = -0.4886189999999999
$ws.Range("D16").Value = "Adding Consistently"
$ws.Range("G16").Value = 1.490929
$ws.Range("H16").Value = 0.137845
$ws.Range("I16").Value = 0.05595799999999995
$ws.Range("D17").Value = "Adding Consistently"
$ws.Range("G17").Value = 1.289392
$ws.Range("H17").Value = 0.1411239999999998
$ws.Range("I17").Value = 0.1583279999999998
$ws.Range("D18").Value = "Reducing Consistently"
$ws.Range("G18").Value = 1.498089
$ws.Range("H18").Value = -0.08724600000000016
$ws.Range("I18").Value = -0.2672380000000001
$ws.Range("D19").Value = "Reducing Consistently"
$ws.Range("G19").Value = 1.374417
$ws.Range("H19").Value = -0.1439440000000001
$ws.Range("I19").Value = -0.1732990000000001
$ws.Range("D20").Value = "Reducing Consistently"
$ws.Range("G20").Value = 1.350438
$ws.Range("H20").Value = -0.06699000000000011
$ws.Range("I20").Value = -0.2641310000000001
$ws.Range("D21").Value = "Reducing Consistently"
$ws.Range("G21").Value = 0.753629
$ws.Range("H21").Value = -0.09561799999999998
$ws.Range("I21").Value = -0.04823199999999994
$ws.Range("D22").Value = "Reducing Consistently"
$ws.Range("G22").Value = 0.689644
$ws.Range("H22").Value = -0.1385580000000001
$ws.Range("I22").Value = -0.1570750000000001
$ws.Range("D23").Value = "Reducing Consistently"
$ws.Range("G23").Value = 0.39986
$ws.Range("H23").Value = -0.03114099999999997
$ws.Range("I23").Value = -0.046537
$ws.Range("D24").Value = "Complete Exit"
$ws.Range("G24").Value = 0.602359
$ws.Range("H24").Value = -0.786067
$ws.Range("I24").Value = -0.602359
$ws.Range("D25").Value = "Complete Exit"
$ws.Range("G25").Value = 9.547155
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = -9.547155
$ws.Range("D26").Value = "Complete Exit"
$ws.Range("G26").Value = 3.947823
$ws.Range("H26").Value = -3.739444
$ws.Range("I26").Value = -3.947823
